$wb = $excel.ActiveWorkbook

$sheetsAno = @(
  "Potencia Acumulada - SIN (MW)",
  "Geracao Periodo Medio (MWMed)",
  "Atendimento a Ponta(MW)",
  "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $sheetsAno) {
  $ws = $wb.Worksheets.Item($name)
  $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Text
  $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Text
  $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Text
  $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Text
}

$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIntervalo.Range("B1").Value = "Intervalo " + $wsIntervalo.Range("B1").Text
$wsIntervalo.Range("C1").Value = "Intervalo " + $wsIntervalo.Range("C1").Text
$wsIntervalo.Range("D1").Value = "Intervalo " + $wsIntervalo.Range("D1").Text
$wsIntervalo.Range("E1").Value = "Intervalo " + $wsIntervalo.Range("E1").Text

$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano " + $wsCusto.Range("B1").Text
